$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line contingencies (line7, line8) are inserted right after line6 (row 7),
# pushing the existing extr1..extr8 rows (formerly rows 8-15) down by two rows
# (to rows 10-17). The name labels are updated accordingly and the associated
# from_bus/to_bus/in_service data for extr1..extr8 move down with them, while
# line7/line8 receive brand new data.

# Row 8 -> line7 (new)
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9 -> line8 (new)
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Row 10 -> extr1 (shifted down from old row 8)
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11 -> extr2 (shifted down from old row 9)
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12 -> extr3 (shifted down from old row 10)
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

# Row 13 -> extr4 (shifted down from old row 11)
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# Row 14 -> extr5 (shifted down from old row 12)
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

# Row 15 -> extr6 (shifted down from old row 13)
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# Row 16 -> extr7 (shifted down from old row 14) - new row, needs style like column A
$c = $ws.Cells.Item(16, 1)
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Font.Bold = $true
$c.Borders.LineStyle = 1
$c.Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# Row 17 -> extr8 (shifted down from old row 15) - new row, needs style like column A
$c = $ws.Cells.Item(17, 1)
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Font.Bold = $true
$c.Borders.LineStyle = 1
$c.Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true

Write-Output "edit complete"
